$d = $word.ActiveDocument

$d.Content.Find.Execute('${roleInspektur}', $true, $false, $false, $false, $false,
                         $true, 1, $false, '(Masukkan jabatan penandatangan)', 2)

$d.Content.Find.Execute('${inspektur}', $true, $false, $false, $false, $false,
                         $true, 1, $false, '(Masukkan nama penandatangan)', 2)
